$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gained two new data rows for a "Femacal de La Calera" Naranja
# Valencia record dated 45021 (2023-04-05). They were inserted just above
# the existing row 1150, pushing the old rows 1150-1232 down to 1152-1234.
$ws.Rows("1150:1151").Insert()

# Seed the two new rows with the row directly below them (which, after the
# insert, still holds the original row-1150/1151 content verbatim) so every
# column that does NOT change (A, B, C, E, F, G, H, I, J, Q, R, T) is
# correct, then overwrite only the cells that actually differ.
$ws.Range("A1152:T1152").Copy($ws.Range("A1150:T1150"))
$ws.Range("A1153:T1153").Copy($ws.Range("A1151:T1151"))

# New row 1150: Valencia / Primera
$ws.Range("D1150").Value = 45021
$ws.Range("K1150").Value = "Valencia"
$ws.Range("L1150").Value = "Primera"
$ws.Range("M1150").Value = 65
$ws.Range("N1150").Value = 8000
$ws.Range("O1150").Value = 8000
$ws.Range("P1150").Value = 8000
$ws.Range("S1150").Value = 615

# New row 1151: Valencia / Segunda
$ws.Range("D1151").Value = 45021
$ws.Range("K1151").Value = "Valencia"
$ws.Range("L1151").Value = "Segunda"
$ws.Range("M1151").Value = 78
$ws.Range("N1151").Value = 7000
$ws.Range("O1151").Value = 7000
$ws.Range("P1151").Value = 7000
$ws.Range("S1151").Value = 538
